$d = $word.ActiveDocument
$paras = $d.Paragraphs

$rsquo = [char]0x2019

$oldAbstract = "For the first lab, students were tasked with introducing themselves to Cypress" + $rsquo + " Programmable System on a Chip (PSoC). This goal was achieved through video tutorials, and guides provided on Cypress" + $rsquo + " official website. By being able to utilize PSoC, the students will be able to complete the rest of the labs in ELC 363." + "`r"
$newAbstract = "For the first lab, students were tasked with introducing themselves to Cypress" + $rsquo + " Programmable System on a Chip (PSoC). This goal was achieved through video tutorials, and guides provided on Cypress" + $rsquo + " official website. By being able to utilize PSoC, the students will be able to complete the rest of the labs in ELC 343."

for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text

    if ($t -eq "Laboratory #2 `r") {
        $p.Range.Text = "Laboratory #1 "
    }
    elseif ($t -eq "Introduction to PSoc Creator`r") {
        $p.Range.Text = "Introduction to PSoC Creator"
    }
    elseif ($t -eq "ELC 363-L2`r") {
        $p.Range.Text = "ELC 343-L2"
    }
    elseif ($t -eq $oldAbstract) {
        $p.Range.Text = $newAbstract
    }
}
